$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the PO numbers in rows 2 and 3
$ws.Range("C2").Value = 1774009
$ws.Range("C3").Value = 1774010

# Row 4's data is no longer needed - clear its contents (keeps border, resets font/fill to plain)
$ws.Range("A4:C4").ClearContents()
$ws.Range("A1").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats -> plain bordered/centered style (s=1)

# Rows 5-7 (B:C) had the bold/shaded "template" formatting left on them - reset to plain style
$ws.Range("A5").Copy()
$ws.Range("B5:C7").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Remove the now-unused trailing blank rows 21-25
$ws.Rows("21:25").Delete()

# Reselect rows 2:3 (matches the final saved selection)
$null = $ws.Rows("2:3").Select()
